$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 87 (shifts old rows 87..221 down to 88..222,
# and Excel's native Insert carries the date-format style of D87 down with it).
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new weekly price observation.
$ws.Cells.Item(87, 1).Value  = 5
$ws.Cells.Item(87, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(87, 3).Value  = "Maule"
$ws.Cells.Item(87, 4).Value  = 44792
$ws.Cells.Item(87, 5).Value  = 7
$ws.Cells.Item(87, 6).Value  = 100112017
$ws.Cells.Item(87, 7).Value  = "Apio"
$ws.Cells.Item(87, 8).Value  = "Americana (o)"
$ws.Cells.Item(87, 9).Value  = "Primera"
$ws.Cells.Item(87, 10).Value = 500
$ws.Cells.Item(87, 11).Value = 10000
$ws.Cells.Item(87, 12).Value = 10000
$ws.Cells.Item(87, 13).Value = 10000
$ws.Cells.Item(87, 14).Value = "`$/docena de matas"
$ws.Cells.Item(87, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(87, 16).Value = 1667
$ws.Cells.Item(87, 17).Value = 6
$ws.Cells.Item(87, 18).Value = "Hortaliza"
